# Update the lattice-multiplication exercise table: replace the 15
# "A x B" problems (and their derived lattice-grid label rows) with a
# new set of problems, cell-for-cell, preserving all formatting.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# New problems, in row-major reading order (row 1: col1,col2,col3; row 2: ...)
$newProblems = @(
    "97 x 69", "78 x 30", "27 x 63",
    "60 x 93", "76 x 12", "97 x 29",
    "93 x 11", "79 x 26", "68 x 37",
    "62 x 10", "54 x 66", "96 x 48",
    "84 x 41", "47 x 21", "38 x 84"
)

$vbreak = [char]11   # manual line break, corresponds to <w:br/>

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$index = 0

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $problem = $newProblems[$index]
        $index = $index + 1

        $parts = $problem.Split("x")
        $a = $parts[0].Trim()
        $b = $parts[1].Trim()

        $aDigits = $a.ToCharArray()
        $bDigits = $b.ToCharArray()
        $a0 = [string]$aDigits[0]
        $a1 = [string]$aDigits[1]
        $b0 = [string]$bDigits[0]
        $b1 = [string]$bDigits[1]

        $line1 = $problem
        $line2 = "  $b0    $b1"
        $line3 = "  ----"
        $line4 = "$a0|    |"
        $line5 = "$a1|    |"

        $newText = "$line1$vbreak$line2$vbreak$line3$vbreak$line4$vbreak$line5"

        $cellRange = $t.Cell($r, $c).Range
        $cellRange.End = $cellRange.End - 1
        $cellRange.Text = $newText
    }
}
